# The sheet previously had a literal "Key"/"Value" header row backing the
# Table2 ListObject, followed by 12 audience rows. The table is converted to
# a headerless table (Table Design > Header Row unchecked) and the leftover
# "Key"/"Value" row is removed from the worksheet, shifting every row (data
# rows and the spacer "E" cells below the table) up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Table Design ribbon -> uncheck "Header Row": the table no longer treats
# row 1 as a header (headerRowCount becomes 0) and its autofilter goes away,
# while the column names ("Key"/"Value") are retained internally.
$lo.ShowHeaders = $false

# The now-redundant "Key"/"Value" row is deleted outright, so the table
# (and all the spacer cells in column E beneath it) shift up by one row.
$ws.Rows(1).Delete()

# Leave the same cell selected/active as in the saved workbook.
$ws.Range("C4").Select() | Out-Null
